# Automatic update of files.
#
# The source data for this "Artfynd" export was re-sorted upstream, which
# swaps the content of rows 13/14 and rows 18/19 (the record identifiers,
# species fields and coordinates move to the other row while the rest of
# each row's cells stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-CellValues($ws, $addr1, $addr2) {
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

# Columns that differ between row 13 and row 14 in the source diff.
$cols1314 = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")
foreach ($col in $cols1314) {
    Swap-CellValues $ws "$col`13" "$col`14"
}

# Columns that differ between row 18 and row 19 in the source diff.
$cols1819 = @("A", "P", "Q", "R", "S")
foreach ($col in $cols1819) {
    Swap-CellValues $ws "$col`18" "$col`19"
}
